$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1, matching the format of the existing
# header cells (bold, bordered, centered) by copying the format from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the "Save" column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
